# Daily GitHub-Actions refresh of the cryptos price/volume table.
# Mirrors the scraped coinranking.com values: updates Price (D) and
# Volume(1h) (E) for each ranked coin, and (for this run) swaps the
# Litecoin / NEARProtocol rows back into rank order (rows 26-27).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. '205.61') need a
# leading apostrophe so Excel stores them as text, same as the original
# inlineStr cells, instead of silently converting them to numeric values.

$ws.Range("D2").Value = '79.849.37'
$ws.Range("E2").Value = '  +4.42%  '
$ws.Range("D3").Value = '3.205.88'
$ws.Range("E3").Value = '  +5.11%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''205.61'
$ws.Range("E5").Value = '  +1.61%  '
$ws.Range("D6").Value = '''634.55'
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.240'
$ws.Range("E8").Value = '  +14.57%  '
$ws.Range("D9").Value = '''0.585'
$ws.Range("E9").Value = '  +5.92%  '
$ws.Range("D10").Value = '3.203.65'
$ws.Range("E10").Value = '  +5.08%  '
$ws.Range("D11").Value = '''0.583'
$ws.Range("E11").Value = '  +33.36%  '
$ws.Range("E12").Value = '  +2.99%  '
$ws.Range("E13").Value = '  +7.05%  '
$ws.Range("E14").Value = '  +19.90%  '
$ws.Range("D15").Value = '3.793.80'
$ws.Range("E15").Value = '  +5.00%  '
$ws.Range("D16").Value = '''32.05'
$ws.Range("E16").Value = '  +8.84%  '
$ws.Range("D17").Value = '79.594.77'
$ws.Range("E17").Value = '  +4.14%  '
$ws.Range("D18").Value = '3.196.97'
$ws.Range("E18").Value = '  +4.61%  '
$ws.Range("D19").Value = '''14.54'
$ws.Range("E19").Value = '  +7.46%  '
$ws.Range("D20").Value = '''3.01'
$ws.Range("E20").Value = '  +29.73%  '
$ws.Range("D21").Value = '''9.23'
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("D22").Value = '''430.64'
$ws.Range("E22").Value = '  +14.87%  '
$ws.Range("D23").Value = '''5.12'
$ws.Range("E23").Value = '  +17.62%  '
$ws.Range("D24").Value = '''11.28'
$ws.Range("E24").Value = '  +13.40%  '
$ws.Range("D25").Value = '3.368.25'
$ws.Range("E25").Value = '  +5.03%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '''77.19'
$ws.Range("E26").Value = '  +4.90%  '
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D27").Value = '''4.76'
$ws.Range("E27").Value = '  +7.55%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  +7.65%  '
$ws.Range("D30").Value = '''9.06'
$ws.Range("E30").Value = '  +8.92%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("E32").Value = '  +5.38%  '
$ws.Range("D33").Value = '''528.24'
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("E34").Value = '  +2.46%  '
$ws.Range("D35").Value = '''0.144'
$ws.Range("E35").Value = '  +27.54%  '
$ws.Range("D36").Value = '''22.99'
$ws.Range("E36").Value = '  +10.10%  '
$ws.Range("E37").Value = '  +11.68%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  +5.70%  '
$ws.Range("D40").Value = '''165.46'
$ws.Range("E40").Value = '  +1.45%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '''193.21'
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '''5.56'
$ws.Range("E44").Value = '  +6.69%  '
$ws.Range("D45").Value = '''0.825'
$ws.Range("E45").Value = '  +2.48%  '
$ws.Range("E46").Value = '  +7.70%  '
$ws.Range("D47").Value = '''1.32'
$ws.Range("E47").Value = '  +4.11%  '
$ws.Range("D48").Value = '''43.31'
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("D49").Value = '''26.01'
$ws.Range("E49").Value = '  +15.51%  '
$ws.Range("D50").Value = '''0.640'
$ws.Range("E50").Value = '  +4.63%  '
$ws.Range("D51").Value = '''2.52'
$ws.Range("E51").Value = '  +1.72%  '
